# Applies the scheduled-runner update to Sheets/Siren_Profits.xlsx
# Updates recomputed leve-profit figures (columns H-N) across all 8 job sheets.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 32
$ws.Range("H32").Value = 7750
$ws.Range("I32").Value = 7750
$ws.Range("K32").Value = 7750
$ws.Range("M32").Value = -7424
# row 45
$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()
# row 88
$ws.Range("H88").Value = 587.7778
$ws.Range("I88").Value = 699.5
$ws.Range("J88").Value = 555.8570999999999
$ws.Range("K88").Value = 699.5
$ws.Range("L88").Value = 555.8570999999999
$ws.Range("M88").Value = -293.5
$ws.Range("N88").Value = -1367.8571
# row 91
$ws.Range("H91").Value = 587.7778
$ws.Range("I91").Value = 699.5
$ws.Range("J91").Value = 555.8570999999999
$ws.Range("K91").Value = 699.5
$ws.Range("L91").Value = 555.8570999999999
$ws.Range("M91").Value = 704.5
$ws.Range("N91").Value = -3363.8571
# row 92
$ws.Range("H92").Value = 750.2857
$ws.Range("I92").Value = 750.75
$ws.Range("J92").Value = 749.6667
$ws.Range("K92").Value = 750.75
$ws.Range("L92").Value = 749.6667
$ws.Range("M92").Value = 497.25
$ws.Range("N92").Value = -3245.6667
# row 116
$ws.Range("H116").Value = 2228861.5
$ws.Range("I116").Value = 5559055.5
$ws.Range("J116").Value = 8732.333000000001
$ws.Range("K116").Value = 5559055.5
$ws.Range("L116").Value = 8732.333000000001
$ws.Range("M116").Value = -5555613.5
$ws.Range("N116").Value = -15616.333
# row 132
$ws.Range("H132").Value = 4672.433
$ws.Range("I132").Value = 5230.0435
$ws.Range("K132").Value = 15690.1305
$ws.Range("M132").Value = -13160.1305
# row 138
$ws.Range("H138").Value = 4490.355
$ws.Range("I138").Value = 1265.3334
$ws.Range("J138").Value = 4968.1357
$ws.Range("K138").Value = 3796.0002
$ws.Range("L138").Value = 14904.4071
$ws.Range("M138").Value = 1343.9998
$ws.Range("N138").Value = -25184.4071

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 61
$ws.Range("H61").Value = 15974
$ws.Range("I61").Value = 19965.666
$ws.Range("J61").Value = 3999
$ws.Range("K61").Value = 19965.666
$ws.Range("L61").Value = 3999
$ws.Range("M61").Value = -19753.666
$ws.Range("N61").Value = -4423
# row 76
$ws.Range("H76").Value = 555555
$ws.Range("J76").Value = 555555
$ws.Range("L76").Value = 555555
$ws.Range("N76").Value = -556231
# row 79
$ws.Range("H79").Value = 555555
$ws.Range("J79").Value = 555555
$ws.Range("L79").Value = 555555
$ws.Range("N79").Value = -557895
# row 110
$ws.Range("H110").Value = 1876.6666
$ws.Range("I110").Value = 1876.6666
$ws.Range("K110").Value = 1876.6666
$ws.Range("M110").Value = 168.3334
# row 136
$ws.Range("H136").Value = 15974
$ws.Range("I136").Value = 19965.666
$ws.Range("J136").Value = 3999
$ws.Range("K136").Value = 59896.99800000001
$ws.Range("L136").Value = 11997
$ws.Range("M136").Value = -57346.99800000001
$ws.Range("N136").Value = -17097
# row 138
$ws.Range("H138").Value = 81184.664
$ws.Range("J138").Value = 81184.664
$ws.Range("L138").Value = 81184.664
$ws.Range("N138").Value = -91464.664

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
# row 25
$ws.Range("H25").Value = 6145.8
$ws.Range("I25").Value = 3576.3333
$ws.Range("K25").Value = 3576.3333
$ws.Range("M25").Value = -3341.3333
# row 134
$ws.Range("H134").Value = 2397.8276
$ws.Range("I134").Value = 1944.9166
$ws.Range("K134").Value = 5834.7498
$ws.Range("M134").Value = -3299.7498

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row 16
$ws.Range("H16").Value = 3121.25
$ws.Range("I16").Value = 3272.7778
$ws.Range("K16").Value = 3272.7778
$ws.Range("M16").Value = -2985.7778
# row 31
$ws.Range("H31").Value = 4763.3213
$ws.Range("I31").Value = 3668.6
$ws.Range("J31").Value = 7500.125
$ws.Range("K31").Value = 3668.6
$ws.Range("L31").Value = 7500.125
$ws.Range("M31").Value = -3373.6
$ws.Range("N31").Value = -8090.125
# row 34
$ws.Range("H34").Value = 4763.3213
$ws.Range("I34").Value = 3668.6
$ws.Range("J34").Value = 7500.125
$ws.Range("K34").Value = 3668.6
$ws.Range("L34").Value = 7500.125
$ws.Range("M34").Value = -3466.6
$ws.Range("N34").Value = -7904.125
# row 86
$ws.Range("H86").Value = 8724.741
$ws.Range("I86").Value = 7595.1665
$ws.Range("K86").Value = 7595.1665
$ws.Range("M86").Value = -6472.1665
# row 89
$ws.Range("H89").Value = 8724.741
$ws.Range("I89").Value = 7595.1665
$ws.Range("K89").Value = 37975.8325
$ws.Range("M89").Value = -32359.8325
# row 113
$ws.Range("H113").Value = 3121.25
$ws.Range("I113").Value = 3272.7778
$ws.Range("K113").Value = 3272.7778
$ws.Range("M113").Value = -1102.7778
# row 141
$ws.Range("H141").Value = 124859.25
$ws.Range("I141").Value = 49500
$ws.Range("J141").Value = 149979
$ws.Range("K141").Value = 49500
$ws.Range("L141").Value = 149979
$ws.Range("M141").Value = -44320
$ws.Range("N141").Value = -160339

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
# row 26
$ws.Range("H26").Value = 194.5
$ws.Range("I26").Value = 114.44444
$ws.Range("K26").Value = 343.33332
$ws.Range("M26").Value = -55.33332000000001
# row 68
$ws.Range("H68").Value = 9559.0625
$ws.Range("I68").Value = 100
$ws.Range("K68").Value = 300
$ws.Range("M68").Value = 511
# row 71
$ws.Range("H71").Value = 9559.0625
$ws.Range("I71").Value = 100
$ws.Range("K71").Value = 900
$ws.Range("M71").Value = 3156
# row 96
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("M96").ClearContents()
# row 100
$ws.Range("H100").Value = 3015
$ws.Range("J100").Value = 2030
$ws.Range("L100").Value = 6090
$ws.Range("N100").Value = -7712
# row 103
$ws.Range("H103").Value = 5165.364
$ws.Range("I103").Value = 6197.3335
$ws.Range("J103").Value = 521.5
$ws.Range("K103").Value = 18592.0005
$ws.Range("L103").Value = 1564.5
$ws.Range("M103").Value = -17713.0005
$ws.Range("N103").Value = -3322.5
# row 131
$ws.Range("H131").Value = 71430370
$ws.Range("I131").Value = 1000000000
$ws.Range("J131").Value = 1932
$ws.Range("K131").Value = 3000000000
$ws.Range("L131").Value = 5796
$ws.Range("M131").Value = -2999994960
$ws.Range("N131").Value = -15876

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
# row 80
$ws.Range("H80").Value = 3141.875
$ws.Range("J80").Value = 3750
$ws.Range("L80").Value = 3750
$ws.Range("N80").Value = -5746
# row 83
$ws.Range("H83").Value = 3141.875
$ws.Range("J83").Value = 3750
$ws.Range("L83").Value = 18750
$ws.Range("N83").Value = -28734
# row 113
$ws.Range("H113").Value = 3899.8
$ws.Range("J113").Value = 4666.3335
$ws.Range("L113").Value = 4666.3335
$ws.Range("N113").Value = -9006.333500000001
# row 132
$ws.Range("H132").Value = 2538.9
$ws.Range("I132").Value = 1547.2727
$ws.Range("J132").Value = 5265.875
$ws.Range("K132").Value = 4641.8181
$ws.Range("L132").Value = 15797.625
$ws.Range("M132").Value = -2111.8181
$ws.Range("N132").Value = -20857.625
# row 134
$ws.Range("H134").Value = 23248.5
$ws.Range("J134").Value = 23248.5
$ws.Range("L134").Value = 69745.5
$ws.Range("N134").Value = -74815.5

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row 22
$ws.Range("H22").Value = 1593.8667
$ws.Range("I22").Value = 2086.2
$ws.Range("J22").Value = 609.2
$ws.Range("K22").Value = 2086.2
$ws.Range("L22").Value = 609.2
$ws.Range("M22").Value = -1791.2
$ws.Range("N22").Value = -1199.2
# row 27
$ws.Range("H27").Value = 1593.8667
$ws.Range("I27").Value = 2086.2
$ws.Range("J27").Value = 609.2
$ws.Range("K27").Value = 2086.2
$ws.Range("L27").Value = 609.2
$ws.Range("M27").Value = -1979.2
$ws.Range("N27").Value = -823.2
# row 46
$ws.Range("H46").Value = 3146.476
$ws.Range("I46").Value = 1298.3334
$ws.Range("K46").Value = 1298.3334
$ws.Range("M46").Value = -1110.3334
# row 68
$ws.Range("H68").Value = 3841.8333
$ws.Range("I68").Value = 3112.8333
$ws.Range("J68").Value = 5299.8335
$ws.Range("K68").Value = 3112.8333
$ws.Range("L68").Value = 5299.8335
$ws.Range("M68").Value = -2363.8333
$ws.Range("N68").Value = -6797.8335
# row 71
$ws.Range("H71").Value = 3841.8333
$ws.Range("I71").Value = 3112.8333
$ws.Range("J71").Value = 5299.8335
$ws.Range("K71").Value = 15564.1665
$ws.Range("L71").Value = 26499.1675
$ws.Range("M71").Value = -11820.1665
$ws.Range("N71").Value = -33987.1675
# row 74
$ws.Range("H74").Value = 69632.5
$ws.Range("I74").Value = 69632.5
$ws.Range("K74").Value = 69632.5
$ws.Range("M74").Value = -68634.5
# row 77
$ws.Range("H77").Value = 69632.5
$ws.Range("I77").Value = 69632.5
$ws.Range("K77").Value = 208897.5
$ws.Range("M77").Value = -203905.5
# row 100
$ws.Range("H100").Value = 2584.7896
$ws.Range("J100").Value = 2656.4167
$ws.Range("L100").Value = 2656.4167
$ws.Range("N100").Value = -3738.4167
# row 127
$ws.Range("H127").Value = 143123660
$ws.Range("J127").Value = 347500
$ws.Range("L127").Value = 347500
$ws.Range("N127").Value = -357420
# row 132
$ws.Range("H132").Value = 3238614.5
$ws.Range("I132").Value = 3930088.8
$ws.Range("K132").Value = 11790266.4
$ws.Range("M132").Value = -11787736.4

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
# row 75
$ws.Range("H75").Value = 36500
$ws.Range("J75").Value = 36500
$ws.Range("L75").Value = 36500
$ws.Range("N75").Value = -38372
# row 78
$ws.Range("H78").Value = 36500
$ws.Range("J78").Value = 36500
$ws.Range("L78").Value = 109500
$ws.Range("N78").Value = -118860

